# daily auto push: 2026-01-21 18:57 UTC
#
# Two new 3-hourly log entries for 2026/01/21 (水) and 2026/01/22 (木) were
# appended to the existing run; because the sheet is sorted chronologically
# they land right before the current row 688 ("2026/12/29"), pushing every
# row from 688 onward down by two (old 688 -> new 690, ..., old 729 -> new 731).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 688; everything at/after 688 shifts down by 2,
# and the sheet's $dimension is updated automatically (A1:D729 -> A1:D731).
$ws.Rows.Item(688).Insert()
$ws.Rows.Item(688).Insert()

# The date column stores text like "2026/01/21", not a real date serial -
# force text formatting before writing so Excel doesn't auto-convert the
# literal into a date value.
$ws.Range("A688:A689").NumberFormat = "@"

$ws.Range("A688").Value = "2026/01/21"
$ws.Range("B688").Value = "水"
$ws.Range("C688").Value = 22
$ws.Range("D688").Value = 201

$ws.Range("A689").Value = "2026/01/22"
$ws.Range("B689").Value = "木"
$ws.Range("C689").Value = 2
$ws.Range("D689").Value = 201

# Drop the temporary text-number-format override so the new cells match the
# unformatted (default-style) look of every other data row.
$ws.Range("A688:A689").ClearFormats()
